$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '65.026.71'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +0.22%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.519.07'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +0.11%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.998'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.36%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '595.47'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +0.50%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '134.44'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -0.99%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.516.62'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +0.09%  '

# Row 8
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -0.03%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.495'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +0.76%  '

# Row 10
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +1.00%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '7.14'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +4.68%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.385'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +0.30%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.104.44'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -0.21%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '27.21'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +0.58%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.0000182'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +0.70%  '

# Row 16
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -0.03%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.510.39'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -0.48%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '64.225.88'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -1.12%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '9.78'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -2.31%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.44'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +2.29%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.71'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -1.66%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '390.31'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +0.61%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.578'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +1.63%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '3.658.16'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -0.21%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '74.30'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +0.73%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.999'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -0.09%  '

# Row 27
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +2.41%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.64'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +20.99%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.79'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +2.10%  '

# Row 30
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +0.58%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.28'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +1.70%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '8.42'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +3.10%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.522.76'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -0.26%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '24.12'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +1.86%  '

# Row 35
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +0.05%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.146'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +2.35%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.28'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +7.16%  '

# Row 38
$ws.Range("B38").NumberFormat = "@"
$ws.Range("B38").Value = 'ImmutableX'
$ws.Range("C38").NumberFormat = "@"
$ws.Range("C38").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.58'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +1.30%  '

# Row 39
$ws.Range("B39").NumberFormat = "@"
$ws.Range("B39").Value = 'Monero'
$ws.Range("C39").NumberFormat = "@"
$ws.Range("C39").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '169.85'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +0.50%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '6.85'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +0.24%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0828'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +4.17%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.822'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +0.88%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '42.62'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +0.80%  '

# Row 44
$ws.Range("B44").NumberFormat = "@"
$ws.Range("B44").Value = 'ONDO'
$ws.Range("C44").NumberFormat = "@"
$ws.Range("C44").Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.24'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +4.15%  '

# Row 45
$ws.Range("B45").NumberFormat = "@"
$ws.Range("B45").Value = 'FirstDigitalUSD'
$ws.Range("C45").NumberFormat = "@"
$ws.Range("C45").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.997'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -0.62%  '

# Row 46
$ws.Range("B46").NumberFormat = "@"
$ws.Range("B46").Value = 'EnergySwap'
$ws.Range("C46").NumberFormat = "@"
$ws.Range("C46").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '25.20'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -3.51%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '4.44'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +1.07%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.66'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +0.37%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '6.93'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +1.93%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.364.18'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -0.78%  '

# Row 51
$ws.Range("B51").NumberFormat = "@"
$ws.Range("B51").Value = 'SuiNetwork'
$ws.Range("C51").NumberFormat = "@"
$ws.Range("C51").Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.895'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +7.29%  '
